# feat: add 2022-Q4 data
#
# Original workbook has two sheets:
#   "总计"     - summary sheet, one data row for 2022-Q3
#   "2022-Q3"  - quarterly detail sheet, one fund row
#
# Target workbook has three sheets:
#   "总计"     - summary sheet, now with a 2022-Q4 row (new) followed by the
#                original 2022-Q3 row
#   "2022-Q4"  - NEW quarterly detail data (3 funds), reusing the
#                "2022-Q3" sheet's tab position/rId
#   "2022-Q3"  - the ORIGINAL quarterly detail sheet content, preserved as a
#                new sheet placed right after "2022-Q4"

$wb = $excel.ActiveWorkbook
$wsTotal = $wb.Worksheets.Item("总计")
$wsQ3 = $wb.Worksheets.Item("2022-Q3")

# ---------------------------------------------------------------------
# 1) Duplicate the existing "2022-Q3" sheet *before* touching its data, so
#    the duplicate keeps the original single-fund content. Then rename the
#    original to "2022-Q4" (it will hold the new data) and the duplicate
#    back to "2022-Q3" (it keeps the old data, now positioned after Q4).
# ---------------------------------------------------------------------
$wsQ3.Copy($null, $wsQ3)
$wsDup = $wb.Worksheets.Item("2022-Q3 (2)")
$wsDup.Name = "2022-Q3-tmp"
$wsQ3.Name = "2022-Q4"
$wsDup.Name = "2022-Q3"

$wsQ4 = $wb.Worksheets.Item("2022-Q4")
$wsQ3New = $wb.Worksheets.Item("2022-Q3")

# ---------------------------------------------------------------------
# 2) Update "总计": shift the existing 2022-Q3 row down to row 3, add the
#    new 2022-Q4 row in row 2 (row 1 is the header).
# ---------------------------------------------------------------------
$wsTotal.Range("A2:D2").Copy()
$wsTotal.Range("A3:D3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("C3").Value = 1
$wsTotal.Range("D3").Value = 0.01

$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 3
$wsTotal.Range("D2").Value = 0.02

# ---------------------------------------------------------------------
# 3) Replace "2022-Q4" sheet's data with the new quarter's figures.
#    Clear the old single-fund row, restyle header + column A to match the
#    "总计" header style, then write the three new fund rows.
# ---------------------------------------------------------------------
$wsQ4.Rows("2:4").ClearContents()

$wsTotal.Range("B1:D1").Copy()
$wsQ4.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$wsTotal.Range("A2").Copy()
$wsQ4.Range("A2:A4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Numeric column A (row index)
$wsQ4.Range("A2").Value = 0
$wsQ4.Range("A3").Value = 1
$wsQ4.Range("A4").Value = 2

# Text columns B:G must stay text (fund codes / percentages as strings),
# so force a text format before writing, then drop back to the default
# "Normal" style (keeps values textual without leaving a numFmt override).
$txtRange = $wsQ4.Range("B2:G4")
$txtRange.NumberFormat = "@"

$wsQ4.Range("B2").Value = "161620"
$wsQ4.Range("C2").Value = "融通核心价值混合（QDII）A"
$wsQ4.Range("D2").Value = "0.51"
$wsQ4.Range("E2").Value = "65.37"
$wsQ4.Range("F2").Value = "2.90"
$wsQ4.Range("G2").Value = "0.0148"

$wsQ4.Range("B3").Value = "005269"
$wsQ4.Range("C3").Value = "华泰柏瑞港股通量化灵活配置混合"
$wsQ4.Range("D3").Value = "0.54"
$wsQ4.Range("E3").Value = "80.96"
$wsQ4.Range("F3").Value = "1.68"
$wsQ4.Range("G3").Value = "0.0091"

$wsQ4.Range("B4").Value = "014127"
$wsQ4.Range("C4").Value = "融通核心价值混合（QDII）C"
$wsQ4.Range("D4").Value = "0.02"
$wsQ4.Range("E4").Value = "65.37"
$wsQ4.Range("F4").Value = "2.90"
$wsQ4.Range("G4").Value = "0.0006"

$txtRange.Style = "Normal"

# Numeric rank column H
$wsQ4.Range("H2").Value = 7
$wsQ4.Range("H3").Value = 9
$wsQ4.Range("H4").Value = 7


